$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 223, shifting the existing rows 223-264 down to 224-265.
$ws.Rows.Item(223).Insert()

# Populate the newly inserted row 223 with the new record.
$ws.Range("A223").Value = 3
$ws.Range("B223").Value = "Femacal de La Calera"
$ws.Range("C223").Value = "Coquimbo"
$ws.Range("D223").Value = 44529
$ws.Range("E223").Value = 5
$ws.Range("F223").Value = 100112040
$ws.Range("G223").Value = "Cilantro"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 120
$ws.Range("K223").Value = 5000
$ws.Range("L223").Value = 5000
$ws.Range("M223").Value = 5000
$ws.Range("N223").Value = "$/docena de atados (3 kilos)"
$ws.Range("O223").Value = "Provincia de Quillota"
$ws.Range("P223").Value = 1667
$ws.Range("Q223").Value = 3
$ws.Range("R223").Value = "Hortaliza"

# Keep the D223 cell formatted the same way as the other date cells in column D.
$ws.Range("D223").NumberFormat = $ws.Range("D224").NumberFormat()
